$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 5 through 28 (data rows beyond the new truncated range)
$ws.Range("A5:A28").EntireRow.Delete() | Out-Null

# Update the remaining data values
$ws.Range("A2").Value = 442
$ws.Range("A3").Value = 523
$ws.Range("A4").Value = 646

# Add an (empty) date-formatted cell at N2 (builtin numFmtId 14)
$ws.Range("N2").NumberFormat = "mm-dd-yy"

# Move the active selection to B2
$ws.Range("B2").Select() | Out-Null
